$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the year header row (row 4) to 2015..2021 ---
$years = 2015,2016,2017,2018,2019,2020,2021
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(4, 4 + $i).Value = $years[$i]
}

# --- Re-point the formatting of row 5 (D5:G5 currently use the custom "164" style,
#     the rest of the row (H5:J5) already uses the plain style (s=7); make the whole
#     row consistent with that plain style before writing new values ---
$ws.Range("H5").Copy() | Out-Null
$ws.Range("D5:G5").PasteSpecial(-4122) | Out-Null

$row5 = 2.2197193775563164,2.1235271668715399,2.7818537161298167,6.7272960584548969,5.1525830614767187,4.4774536255935971,4.6024666695867751
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 4 + $i).Value = $row5[$i]
}

# --- Row 6: D6 keeps its style, E6:J6 should use the style that M6:P6 used ---
$ws.Range("M6").Copy() | Out-Null
$ws.Range("E6:J6").PasteSpecial(-4122) | Out-Null

$row6 = 2.2322863217945752,2.8603553109638966,3.113207036164539,6.2970593463100784,4.8617746111834492,2.6715092780025032,4.3694509108608912
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 4 + $i).Value = $row6[$i]
}

$excel.CutCopyMode = 0

# --- Remove the now-unneeded trailing year columns (K:P), shrinking the table to A:J ---
$ws.Range("K1:P1").EntireColumn.Delete() | Out-Null

# --- Column widths: D:J become a uniform custom width ---
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 8.6

# --- Selection / active cell moves ---
$ws.Range("K16").Select() | Out-Null
